# Update countries & provincias Spain
# Refreshes the "Pais" ranking sheet to a later data pull (24 Abril 2020,
# 10:22 instead of 09:52). A handful of countries' rows changed position
# in the (descending, by "Casos totales") ranking and/or picked up new
# counts; a few other rows only got their active/recovered counts revised.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 10:22"

# --- Austria (row 24): Casos activos / Recuperados revised ------------
$ws.Cells.Item(24, 4).Value = 11872
$ws.Cells.Item(24, 5).Value = 2608

# --- Rows 44/45: Chequia and Filipinas swap rank, both get new data ---
# Row 44 now holds Filipinas (freshly updated figures).
$ws.Cells.Item(44, 1).Value = "Filipinas"
$ws.Cells.Item(44, 2).Value = 7192
$ws.Cells.Item(44, 3).Value = 211
$ws.Cells.Item(44, 4).Value = 762
$ws.Cells.Item(44, 5).Value = 5953
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(44, 7).Value = 15
$ws.Cells.Item(44, 8).Value = 477

# Row 45 now holds Chequia (its previous figures, just moved down a slot).
$ws.Cells.Item(45, 1).Value = "Chequia"
$ws.Cells.Item(45, 2).Value = 7188
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(45, 4).Value = 2186
$ws.Cells.Item(45, 5).Value = 4789
$ws.Cells.Item(45, 6).Value = 76
$ws.Cells.Item(45, 7).Value = 3
$ws.Cells.Item(45, 8).Value = 213

# --- Moldavia (row 59): Casos activos / Recuperados revised -----------
$ws.Cells.Item(59, 4).Value = 755
$ws.Cells.Item(59, 5).Value = 2091

# --- Uzbekistan (row 68): Casos activos / Recuperados revised ---------
$ws.Cells.Item(68, 4).Value = 592
$ws.Cells.Item(68, 5).Value = 1179

# --- Rows 78/79: Camerun and Eslovaquia swap rank, both get new data --
# Row 78 now holds Eslovaquia (freshly updated figures).
$ws.Cells.Item(78, 1).Value = "Eslovaquia"
$ws.Cells.Item(78, 2).Value = 1360
$ws.Cells.Item(78, 3).Value = 35
$ws.Cells.Item(78, 4).Value = 355
$ws.Cells.Item(78, 5).Value = 988
$ws.Cells.Item(78, 6).Value = 8
$ws.Cells.Item(78, 7).Value = 2
$ws.Cells.Item(78, 8).Value = 17

# Row 79 now holds Camerun (its previous figures, just moved down a slot).
$ws.Cells.Item(79, 1).Value = "Camerun"
$ws.Cells.Item(79, 2).Value = 1334
$ws.Cells.Item(79, 3).Value = 0
$ws.Cells.Item(79, 4).Value = 668
$ws.Cells.Item(79, 5).Value = 623
$ws.Cells.Item(79, 6).Value = 33
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 8).Value = 43

# --- Etiopia (row 140): Casos totales / Nuevos casos / Recuperados ----
$ws.Cells.Item(140, 2).Value = 117
$ws.Cells.Item(140, 3).Value = 1
$ws.Cells.Item(140, 5).Value = 93
